# Weekly update: insert a new Alcachofa price record as the new first data
# row (row 3), pushing the previously-existing rows (old rows 3-33) down by
# one row each (they become rows 4-34).
#
# The new row 3 duplicates the data that was in the old row 3, except for a
# new date (column D), reflecting the new weekly observation added to the
# series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 3:33 down to 4:34 by inserting a new blank row at 3.
$ws.Rows("3:3").Insert()

# Populate the new row 3 with the new weekly record.
$ws.Range("A3").Value = 7
$ws.Range("B3").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C3").Value = "Ñuble"
$ws.Range("D3").Value = 44515
$ws.Range("E3").Value = 16
$ws.Range("F3").Value = 100112013
$ws.Range("G3").Value = "Alcachofa"
$ws.Range("H3").Value = "Madrigal"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 120
$ws.Range("K3").Value = 11000
$ws.Range("L3").Value = 12000
$ws.Range("M3").Value = 11500
$ws.Range("N3").Value = "`$/caja 40 unidades"
$ws.Range("O3").Value = "Provincia del Elquí"
$ws.Range("P3").Value = 288
$ws.Range("Q3").Value = 40
$ws.Range("R3").Value = "Hortaliza"
